{"js": "// The document already contains the \"Overfitting\" / \"under fitting\" Q&A\n// paragraphs. The only substantive content change is inside the\n// \"Overfitting occurs ...\" paragraph: the phrase \"has learning\" is\n// replaced by \"capture\" (formatting/run properties are unchanged).\nconst body = context.document.body;\n\nconst results = body.search(\"has learning\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\n// Replace every match (there is exactly one in this document) in place,\n// preserving the surrounding run formatting (plain, lang=\"en-US\").\nfor (let i = results.items.length - 1; i >= 0; i--) {\n  results.items[i].insertText(\"capture\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# The document already contains the \"Overfitting\" / \"under fitting\" Q&A\n# paragraphs. The only substantive content change is inside the\n# \"Overfitting occurs ...\" paragraph: the phrase \"has learning\" is\n# replaced by \"capture\" (formatting/run properties are unchanged).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute(\n  \"has learning\",   # FindText\n  $true,            # MatchCase\n  $false,           # MatchWholeWord\n  $false,           # MatchWildcards\n  $false,           # MatchSoundsLike\n  $false,           # MatchAllWordForms\n  $true,            # Forward\n  1,                # Wrap (wdFindContinue)\n  $false,           # Format\n  \"capture\",        # ReplaceWith\n  2                 # Replace (wdReplaceAll)\n)\n"}
